$d = $word.ActiveDocument

# 1. Replace the original title text with the new one.
$d.Content.Find.Execute("Prueba1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ESTO ES UNA PRUEBA Y UN MODELO DE EJEMPLO", 2) | Out-Null

# 2. Append a blank paragraph, then 6 repetitions of
#    ("LOREM IPSUM LOREM IPSUM" paragraph + blank paragraph).
for ($i = 0; $i -lt 6; $i++) {
    # blank paragraph
    $last = $d.Paragraphs($d.Paragraphs.Count)
    $last.Range.InsertParagraphAfter()

    # "LOREM IPSUM LOREM IPSUM" paragraph
    $last = $d.Paragraphs($d.Paragraphs.Count)
    $last.Range.InsertParagraphAfter()
    $newLast = $d.Paragraphs($d.Paragraphs.Count)
    $newLast.Range.Text = "LOREM IPSUM LOREM IPSUM"
}

# 3. Final trailing blank paragraph.
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
